$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 7106
$ws.Range("I62").Value = 5994
$ws.Range("K62").Value = 5994
$ws.Range("M62").Value = -5370
$ws.Range("H65").Value = 7106
$ws.Range("I65").Value = 5994
$ws.Range("K65").Value = 29970
$ws.Range("M65").Value = -26850
$ws.Range("H98").Value = 1351.3077
$ws.Range("I98").Value = 842.7
$ws.Range("J98").Value = 3046.6667
$ws.Range("K98").Value = 842.7
$ws.Range("L98").Value = 3046.6667
$ws.Range("M98").Value = 655.3
$ws.Range("N98").Value = -6042.6667
$ws.Range("H112").Value = 1519.8667
$ws.Range("J112").Value = 1583.3077
$ws.Range("L112").Value = 4749.9231
$ws.Range("N112").Value = -6965.9231
$ws.Range("H122").Value = 1351.3077
$ws.Range("I122").Value = 842.7
$ws.Range("J122").Value = 3046.6667
$ws.Range("K122").Value = 2528.1
$ws.Range("L122").Value = 9140.000100000001
$ws.Range("M122").Value = -78.10000000000036
$ws.Range("N122").Value = -14040.0001
$ws.Range("H129").Value = 1682.6666
$ws.Range("I129").Value = 830.4375
$ws.Range("J129").Value = 8500.5
$ws.Range("K129").Value = 2491.3125
$ws.Range("L129").Value = 25501.5
$ws.Range("M129").Value = 2508.6875
$ws.Range("N129").Value = -35501.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 21996.334
$ws.Range("I28").Value = 17735
$ws.Range("K28").Value = 17735
$ws.Range("M28").Value = -17543
$ws.Range("H32").Value = 4239.055
$ws.Range("I32").Value = 3244.0488
$ws.Range("K32").Value = 3244.0488
$ws.Range("M32").Value = -2957.0488
$ws.Range("H97").Value = 2787.1428
$ws.Range("I97").Value = 777.64703
$ws.Range("K97").Value = 777.64703
$ws.Range("M97").Value = -281.64703
$ws.Range("H99").Value = 21996.334
$ws.Range("I99").Value = 17735
$ws.Range("K99").Value = 17735
$ws.Range("M99").Value = -14740
$ws.Range("H124").Value = 51490.812
$ws.Range("J124").Value = 51490.812
$ws.Range("L124").Value = 51490.812
$ws.Range("N124").Value = -61310.812

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 71431150
$ws.Range("I94").Value = 125002000
$ws.Range("K94").Value = 125002000
$ws.Range("M94").Value = -125001549
$ws.Range("H134").Value = 137945
$ws.Range("I134").Value = 151538.9
$ws.Range("K134").Value = 454616.7
$ws.Range("M134").Value = -452081.7

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 28508.385
$ws.Range("I31").Value = 37247.035
$ws.Range("J31").Value = 6264.5454
$ws.Range("K31").Value = 37247.035
$ws.Range("L31").Value = 6264.5454
$ws.Range("M31").Value = -36952.035
$ws.Range("N31").Value = -6854.5454
$ws.Range("H34").Value = 28508.385
$ws.Range("I34").Value = 37247.035
$ws.Range("J34").Value = 6264.5454
$ws.Range("K34").Value = 37247.035
$ws.Range("L34").Value = 6264.5454
$ws.Range("M34").Value = -37045.035
$ws.Range("N34").Value = -6668.5454
$ws.Range("H93").Value = 14097.286
$ws.Range("I93").Value = 11372.167
$ws.Range("K93").Value = 11372.167
$ws.Range("M93").Value = -9500.166999999999
$ws.Range("H97").Value = 23967.857
$ws.Range("J97").Value = 23967.857
$ws.Range("L97").Value = 23967.857
$ws.Range("N97").Value = -25949.857
$ws.Range("H132").Value = 2985.3157
$ws.Range("I132").Value = 2781
$ws.Range("J132").Value = 3557.4
$ws.Range("K132").Value = 8343
$ws.Range("L132").Value = 10672.2
$ws.Range("M132").Value = -5813
$ws.Range("N132").Value = -15732.2
$ws.Range("H133").Value = 90138
$ws.Range("J133").Value = 90138
$ws.Range("L133").Value = 90138
$ws.Range("N133").Value = -95198

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 437
$ws.Range("I98").Value = 657.75
$ws.Range("K98").Value = 1973.25
$ws.Range("M98").Value = -475.25
$ws.Range("H129").Value = 89818.56
$ws.Range("I129").Value = 286718.72
$ws.Range("J129").Value = 3674.75
$ws.Range("K129").Value = 860156.1599999999
$ws.Range("L129").Value = 11024.25
$ws.Range("M129").Value = -855156.1599999999
$ws.Range("N129").Value = -21024.25
$ws.Range("H131").Value = 1749.9
$ws.Range("J131").Value = 1923.15
$ws.Range("L131").Value = 5769.450000000001
$ws.Range("N131").Value = -15849.45

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4427
$ws.Range("I80").Value = 3333
$ws.Range("K80").Value = 3333
$ws.Range("M80").Value = -2335
$ws.Range("H83").Value = 4427
$ws.Range("I83").Value = 3333
$ws.Range("K83").Value = 16665
$ws.Range("M83").Value = -11673
$ws.Range("H102").Value = 1684.2667
$ws.Range("I102").Value = 1707.862
$ws.Range("K102").Value = 1707.862
$ws.Range("M102").Value = -85.86200000000008
$ws.Range("H132").Value = 4439.2573
$ws.Range("J132").Value = 7015.4546
$ws.Range("L132").Value = 21046.3638
$ws.Range("N132").Value = -26106.3638

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 20495
$ws.Range("J20").Value = 20495
$ws.Range("L20").Value = 20495
$ws.Range("N20").Value = -20947
$ws.Range("H22").Value = 3139.4
$ws.Range("J22").Value = 1099.5
$ws.Range("L22").Value = 1099.5
$ws.Range("N22").Value = -1689.5
$ws.Range("H27").Value = 3139.4
$ws.Range("J27").Value = 1099.5
$ws.Range("L27").Value = 1099.5
$ws.Range("N27").Value = -1313.5
$ws.Range("H40").Value = 10582.049
$ws.Range("I40").Value = 5380.4644
$ws.Range("K40").Value = 5380.4644
$ws.Range("M40").Value = -5244.4644
$ws.Range("H55").Value = 1065.1333
$ws.Range("I55").Value = 1195.1666
$ws.Range("K55").Value = 1195.1666
$ws.Range("M55").Value = -1022.1666
$ws.Range("H61").Value = 15753.36
$ws.Range("I61").Value = 1589.9412
$ws.Range("J61").Value = 45850.625
$ws.Range("K61").Value = 1589.9412
$ws.Range("L61").Value = 45850.625
$ws.Range("M61").Value = -1387.9412
$ws.Range("N61").Value = -46254.625
$ws.Range("H69").Value = 50000
$ws.Range("J69").Value = 50000
$ws.Range("L69").Value = 50000
$ws.Range("N69").Value = -51622
$ws.Range("H72").Value = 50000
$ws.Range("J72").Value = 50000
$ws.Range("L72").Value = 150000
$ws.Range("N72").Value = -158112
$ws.Range("H82").Value = 3942.7058
$ws.Range("I82").Value = 4279
$ws.Range("J82").Value = 3707.3
$ws.Range("K82").Value = 4279
$ws.Range("L82").Value = 3707.3
$ws.Range("M82").Value = -3918
$ws.Range("N82").Value = -4429.3
$ws.Range("H85").Value = 3942.7058
$ws.Range("I85").Value = 4279
$ws.Range("J85").Value = 3707.3
$ws.Range("K85").Value = 4279
$ws.Range("L85").Value = 3707.3
$ws.Range("M85").Value = -3031
$ws.Range("N85").Value = -6203.3
$ws.Range("H96").Value = 29000
$ws.Range("J96").Value = 29000
$ws.Range("L96").Value = 29000
$ws.Range("N96").Value = -34492
$ws.Range("H99").Value = 30285
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H113").Value = 15753.36
$ws.Range("I113").Value = 1589.9412
$ws.Range("J113").Value = 45850.625
$ws.Range("K113").Value = 1589.9412
$ws.Range("L113").Value = 45850.625
$ws.Range("M113").Value = 580.0588
$ws.Range("N113").Value = -50190.625

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1832.1316
$ws.Range("I122").Value = 1898.0385
$ws.Range("K122").Value = 5694.1155
$ws.Range("M122").Value = -3244.1155
$ws.Range("H132").Value = 2796.5
$ws.Range("I132").Value = 2061.6191
$ws.Range("J132").Value = 5883
$ws.Range("K132").Value = 6184.8573
$ws.Range("L132").Value = 17649
$ws.Range("M132").Value = -3654.8573
$ws.Range("N132").Value = -22709
